$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the existing row 186, pushing rows 186-217
# down to 187-218 (weekly data point added for Brócoli - Macroferia
# Regional de Talca).
$ws.Rows.Item(186).Insert()

$ws.Cells.Item(186, 1).Value = 5
$ws.Cells.Item(186, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(186, 3).Value = "Maule"
$ws.Cells.Item(186, 4).Value2 = 44504
$ws.Cells.Item(186, 5).Value = 7
$ws.Cells.Item(186, 6).Value = 100112023
$ws.Cells.Item(186, 7).Value = "Brócoli"
$ws.Cells.Item(186, 8).Value = "Sin especificar"
$ws.Cells.Item(186, 9).Value = "Primera"
$ws.Cells.Item(186, 10).Value = 6000
$ws.Cells.Item(186, 11).Value = 500
$ws.Cells.Item(186, 12).Value = 500
$ws.Cells.Item(186, 13).Value = 500
$ws.Cells.Item(186, 14).Value = "$/unidad"
$ws.Cells.Item(186, 15).Value = "Región del Maule"
$ws.Cells.Item(186, 16).Value = 500
$ws.Cells.Item(186, 17).Value = 1
$ws.Cells.Item(186, 18).Value = "Hortaliza"
